# Remove the blank paragraph, the "Ver no Jupiter..." paragraph and the
# "(c) 2020 ..." footer paragraph that used to follow the
# "LOM3099: Estatica (Requisito)" line, leaving just the single blank
# paragraph that originally sat right before the page-break paragraph.

$d = $word.ActiveDocument

$requisito = "LOM3099: Est" + [char]0xE1 + "tica (Requisito)"
$jupiter   = "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyright = [char]0xA9 + " 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"

# Locate the "LOM3099: Estatica (Requisito)" paragraph and, right after
# it, the "Ver no Jupiter..." / copyright paragraphs, by scanning the
# Paragraphs collection (robust to any Find/Range quirks).
$anchorIndex = -1
$jupiterIndex = -1
$copyrightIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd()
    if ($text -eq $requisito) {
        $anchorIndex = $i
    } elseif ($anchorIndex -gt 0 -and $jupiterIndex -lt 0 -and $text -eq $jupiter) {
        $jupiterIndex = $i
    } elseif ($jupiterIndex -gt 0 -and $copyrightIndex -lt 0 -and $text -eq $copyright) {
        $copyrightIndex = $i
        break
    }
}

if ($anchorIndex -gt 0 -and $jupiterIndex -eq ($anchorIndex + 2) -and $copyrightIndex -eq ($jupiterIndex + 1)) {
    # Delete from the (blank) paragraph right after the anchor through
    # the end of the copyright paragraph (including all paragraph marks
    # in between), which leaves the blank paragraph that originally
    # followed the copyright line now directly after the anchor.
    $firstToRemove = $d.Paragraphs.Item($anchorIndex + 1)
    $lastToRemove = $d.Paragraphs.Item($copyrightIndex)
    $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End).Delete()
}
